$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before the existing "state_total" column (L), shifting
# the old L column (and its data) one to the right into M.
$ws.Columns("L").Insert()

# New column header + data: "March_16", all zero for every state row.
$ws.Range("L1").Value = "March_16"
$ws.Range("L2:L52").Value = 0

# The insert leaves a stray formatted-but-empty cell at L54 (the blank
# totals-row formatting band); remove it so only M54 carries that formatting.
$ws.Range("L54").Clear()

# Update the recorded selection/active cell to match the finished edit.
$ws.Range("O46").Select()
